$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column for "Unit" between Quantity (E) and Price (old F)
$ws.Range("F1").EntireColumn.Insert()

# Insert a new column for "Cost_Text" after Cost (now H)
$ws.Range("I1").EntireColumn.Insert()

# Update header row
$ws.Range("A1").Value = "INDEX"
$ws.Range("B1").Value = "DATE"
$ws.Range("C1").Value = "ORGANIZATION"
$ws.Range("D1").Value = "ITEM"
$ws.Range("E1").Value = "QUANTITY"
$ws.Range("F1").Value = "UNIT"
$ws.Range("G1").Value = "PRICE"
$ws.Range("H1").Value = "COST"
$ws.Range("I1").Value = "COST_TEXT"

# Fill new Unit column values
$ws.Range("F2").Value = "counts"
$ws.Range("F3").Value = "bulbs"

# Fill new Cost_Text column values
$ws.Range("I2").Value = "ยี่สิบห้าบาท"
$ws.Range("I3").Value = "สองบาท"

# Adjust column C width
$ws.Columns.Item(3).ColumnWidth = 13.25

# Update selection to match target
$null = $ws.Range("H30").Select()
